$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "271.47"
    3  = "21.47"
    4  = "6.293"
    5  = "0.06250"
    6  = "3.568"
    7  = "6.590"
    8  = "1.387"
    9  = "0.8305"
    10 = "0.01350"
    11 = "0.1573"
    12 = "0.08273"
    13 = "0.03375"
    14 = "0.03184"
    16 = "3.603"
    17 = "0.001653"
    18 = "0.04568"
    19 = "0.006319"
    20 = "0.005879"
    21 = "0.001064"
    22 = "0.0001481"
    23 = "3.754"
    24 = "2.366"
    25 = "0.3302"
    26 = "0.1246"
    28 = "0.0002717"
    40 = "0.04683"
    41 = "0.006980"
    42 = "0.1170"
    43 = "0.003612"
    44 = "0.01184"
    45 = "0.00006019"
    46 = "0.0009858"
    47 = "0.00000000741"
    48 = "0.7793"
    49 = "0.002388"
    50 = "0.00002372"
    51 = "0.01225"
}

foreach ($row in $values.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $values[$row]
    $cell.Style = "Normal"
}
